# The commit swaps the contents of the two theme parts in the deck:
#   ppt/theme/theme1.xml  "Office Theme" (Office colours)  -> becomes "Integral" (Red Violet colours)
#   ppt/theme/theme2.xml  "Integral" (Red Violet colours)  -> becomes "Office Theme" (Office colours)
#
# The slide master / presentation's live theme (the one reachable from the
# PowerPoint object model as SlideMaster.Theme / NotesMaster.Theme / etc,
# all of which resolve to the same design) is the "Integral" one, so we
# restore it to the stock "Office Theme" colour scheme -- i.e. apply the
# swap to the theme that the object model exposes.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# ThemeColorScheme items are ordered: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink.
# RGB is a VBA-style BGR-packed long (B*65536 + G*256 + R).

$scheme.Item(1).RGB  = 0         # dk1      000000
$scheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$scheme.Item(3).RGB  = 6968388   # dk2      44546A
$scheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$scheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$scheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$scheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$scheme.Item(8).RGB  = 49407     # accent4  FFC000
$scheme.Item(9).RGB  = 12874308  # accent5  4472C4
$scheme.Item(10).RGB = 4697456   # accent6  70AD47
$scheme.Item(11).RGB = 12673797  # hlink    0563C1
$scheme.Item(12).RGB = 7491477   # folHlink 954F72
